# Tue, Mar 31, 2020  3:06:12 PM
#
# 1) Swap the presentation's theme colour scheme from "Integral" to the
#    stock "Office Theme" palette (dk1/lt1 are already identical black/white
#    in both palettes, so only the remaining 10 slots need to change).
# 2) Re-style the table on slide 16 with the built-in table style
#    {FD3B1DF4-1BF7-4570-9F3F-D3CEC60A4145}.

$p = $ppt.ActivePresentation

# --- 1) Theme colour scheme: Integral -> Office -------------------------
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72

# --- 2) Table style on slide 16 ------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{FD3B1DF4-1BF7-4570-9F3F-D3CEC60A4145}")
    }
}
